# Generate Report for Handback
# The 0b259ff6-c5d9-4fbe-96bc-cfa8dc3bbdaf file has now been handed back and is
# in sync with en-US for both the zh-cn and de-de locales. Update the Overview
# sheet and each locale sheet's row for that file accordingly.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: the 0b259ff6 row (row 2) status columns for zh-cn (B) and
# de-de (C) move from "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# zh-cn sheet: the 0b259ff6 row (row 2) gets a Status update, a Latest Target
# File + Latest Handback File hyperlink pair, and a real Latest Handback
# DateTime (replacing the 0001-01-01 00:00:00 placeholder).
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusHandedBack

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/15ac8d4bf61ccc6867425de8c3c831398021fc23/e2e/0b259ff6-c5d9-4fbe-96bc-cfa8dc3bbdaf.md",
    "",
    "",
    "0b259ff6-c5d9-4fbe-96bc-cfa8dc3bbdaf.md"
)

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33693316306daf2024132236c9f4dd42dd800491/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0b259ff6-c5d9-4fbe-96bc-cfa8dc3bbdaf.a3085c5ee713bdc1999638587fcbd6e2c8e3b609.zh-cn.xlf",
    "",
    "",
    "0b259ff6-c5d9-4fbe-96bc-cfa8dc3bbdaf.a3085c5ee713bdc1999638587fcbd6e2c8e3b609.zh-cn.xlf"
)

$wsZhCn.Range("H2").Value = "2016-03-18 07:10:42"

# ---------------------------------------------------------------------------
# de-de sheet: same shape of edit as zh-cn, different target file names / time.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusHandedBack

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/15ac8d4bf61ccc6867425de8c3c831398021fc23/e2e/0b259ff6-c5d9-4fbe-96bc-cfa8dc3bbdaf.md",
    "",
    "",
    "0b259ff6-c5d9-4fbe-96bc-cfa8dc3bbdaf.md"
)

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/996bda6acf7a86aeb48b2a35f8a115dc543f9840/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0b259ff6-c5d9-4fbe-96bc-cfa8dc3bbdaf.a3085c5ee713bdc1999638587fcbd6e2c8e3b609.de-de.xlf",
    "",
    "",
    "0b259ff6-c5d9-4fbe-96bc-cfa8dc3bbdaf.a3085c5ee713bdc1999638587fcbd6e2c8e3b609.de-de.xlf"
)

$wsDeDe.Range("H2").Value = "2016-03-18 07:10:46"
